# Refresh cryptocurrency price/volume snapshot (scheduled GitHub Actions update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.609.39'
$ws.Range('E2').Value = '  +0.72%  '
$ws.Range('D3').Value = '1.639.58'
$ws.Range('E3').Value = '  +0.91%  '
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').Value = '''214.44'
$ws.Range('E5').Value = '  +0.83%  '
$ws.Range('D6').Value = '''0.502'
$ws.Range('E6').Value = '  +1.60%  '
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('D9').Value = '''0.0624'
$ws.Range('E9').Value = '  +0.70%  '
$ws.Range('D10').Value = '''19.03'
$ws.Range('E10').Value = '  +0.59%  '
$ws.Range('E11').Value = '  +0.51%  '
$ws.Range('D12').Value = '1.867.96'
$ws.Range('E12').Value = '  +0.91%  '
$ws.Range('D13').Value = '1.630.29'
$ws.Range('E13').Value = '  -0.19%  '
$ws.Range('D14').Value = '''4.15'
$ws.Range('E14').Value = '  +1.42%  '
$ws.Range('D15').Value = '''0.527'
$ws.Range('E15').Value = '  +1.38%  '
$ws.Range('D16').Value = '''64.75'
$ws.Range('E16').Value = '  +1.29%  '
$ws.Range('D17').Value = '26.619.69'
$ws.Range('E17').Value = '  +0.79%  '
$ws.Range('D18').Value = '0.0₃0740'
$ws.Range('E18').Value = '  +0.27%  '
$ws.Range('D19').Value = '''214.96'
$ws.Range('E19').Value = '  -0.10%  '
$ws.Range('D20').Value = '''1.00'
$ws.Range('E20').Value = '  -0.25%  '
$ws.Range('D21').Value = '''4.33'
$ws.Range('E21').Value = '  +0.85%  '
$ws.Range('D22').Value = '''6.23'
$ws.Range('E22').Value = '  +0.32%  '
$ws.Range('D23').Value = '''9.43'
$ws.Range('E23').Value = '  +1.79%  '
$ws.Range('D24').Value = '''2.20'
$ws.Range('E24').Value = '  +10.97%  '
$ws.Range('D25').Value = '''144.74'
$ws.Range('E25').Value = '  -2.08%  '
$ws.Range('D26').Value = '''1.00'
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('E27').Value = '  -0.51%  '
$ws.Range('D28').Value = '''7.11'
$ws.Range('E28').Value = '  +4.22%  '
$ws.Range('D29').Value = '''15.64'
$ws.Range('E29').Value = '  +0.69%  '
$ws.Range('D30').Value = '''0.0511'
$ws.Range('E30').Value = '  +0.99%  '
$ws.Range('E31').Value = '  +1.21%  '
$ws.Range('D32').Value = '''3.35'
$ws.Range('E32').Value = '  +1.05%  '
$ws.Range('D33').Value = '''2.98'
$ws.Range('E33').Value = '  +1.64%  '
$ws.Range('D34').Value = '1.274.17'
$ws.Range('E34').Value = '  +5.11%  '
$ws.Range('E35').Value = '  +2.57%  '
$ws.Range('E36').Value = '  +0.97%  '
$ws.Range('D37').Value = '''0.0177'
$ws.Range('E37').Value = '  +2.85%  '
$ws.Range('D38').Value = '''0.528'
$ws.Range('E38').Value = '  +5.85%  '
$ws.Range('D39').Value = '''0.823'
$ws.Range('E39').Value = '  +3.59%  '
$ws.Range('D40').Value = '''1.00'
$ws.Range('E40').Value = '  -0.23%  '
$ws.Range('D41').Value = '''0.808'
$ws.Range('E41').Value = '  +2.09%  '
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('D43').Value = '''5.41'
$ws.Range('E43').Value = '  +0.85%  '
$ws.Range('D44').Value = '1.777.99'
$ws.Range('E44').Value = '  +1.01%  '
$ws.Range('D45').Value = '''91.27'
$ws.Range('E45').Value = '  -1.33%  '
$ws.Range('D46').Value = '''58.88'
$ws.Range('E46').Value = '  +7.86%  '
$ws.Range('E47').Value = '  +1.23%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '''0.0514'
$ws.Range('E48').Value = '  +0.74%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '''7.71'
$ws.Range('E49').Value = '  +1.37%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = '''0.0960'
$ws.Range('E50').Value = '  +1.63%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = '''0.405'
$ws.Range('E51').Value = '  -0.67%  '
